# psychopy-2/runs.xlsx edit
# - Insert a new "contextRole" column (C) with values irrelevant / modulatory / additive
#   (randomly shuffled relative to the existing rows)
# - Fix a typo in the restaurant-name string ("Café" -> "Cafe", ASCII only)
# - Insert a new row 2 containing a long merged instructional note (italic, left aligned)
# - Reorder rows underneath so food-prefix / restaurant-name pairings line up with their
#   shuffled context role

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a blank row right below the header for the instructions banner ---
$ws.Range("A2").EntireRow.Insert()

# --- Header row (row 1): add contextRole header in column C ---
$ws.Range("C1").Value = "contextRole"
$ws.Range("C1").Font.Bold = $true

# --- Row 2: merged instructions banner ---
$ws.Range("A2").Value = "Okay here's how this works. Each row is a different run -- a list of restaurant names + the prefix of the filenames of its corresponding foods. The list of restaurant names should be comma-separated and should contain exactly 3 restaurants (no unicode characters please -- stick to ASCII). For each food prefix, there should be 3 corresponding files in the 'foods' directory in the format ""[food prefix][cueId].png"". For example, for the prefix ""mexican_food"", there should be files ""mexican_food0.png"", ""mexican_food1.png"" and ""mexican_food2.png"".  The third column is randomly shuffled so different foods/restaurants can be used with different context roles across subjects. The context role also corresponds to the filename where the context/cue pairs are, in the format ""[contextRole].xlsx"". In each of these files, the contextId and cueId are matched with the restaurants and foods, respectively. Sorry it's so complicated but there is no other way to make it work with psychopy"

# italic (non-bold), left-aligned style for the banner text, set BEFORE merging so the
# merged range starts from the correct top-left formatting
$ws.Range("A2").Font.Bold = $false
$ws.Range("A2").Font.Italic = $true
$ws.Range("A2").HorizontalAlignment = -4131

# bold+italic, left-aligned style, pre-seeded on B2 before merging
$ws.Range("B2").Font.Bold = $true
$ws.Range("B2").Font.Italic = $true
$ws.Range("B2").HorizontalAlignment = -4131

$ws.Range("A2:C2").Merge()

# merging re-applies A2's style to the whole range, so restore B2 & C2's
# bold+italic, left-aligned formatting afterwards
$ws.Range("B2").Font.Bold = $true
$ws.Range("B2").Font.Italic = $true
$ws.Range("B2").HorizontalAlignment = -4131

$ws.Range("C2").Font.Bold = $true
$ws.Range("C2").Font.Italic = $true
$ws.Range("C2").HorizontalAlignment = -4131

# --- Data rows (now 3,4,5) - rewrite restaurant / food pairing and shuffled context role ---
$ws.Range("A3").Value = "Molina's Cantina,Restaurante Arroyo,El Coyote Cafe"
$ws.Range("B3").Value = "mexican_food"
$ws.Range("C3").Value = "irrelevant"

$ws.Range("A4").Value = "Le Parisien,Chez Toinette,Au Petit Sud Ouest"
$ws.Range("B4").Value = "french_food"
$ws.Range("C4").Value = "modulatory"

$ws.Range("A5").Value = "Lau's Dim Sum Bar,OO Kook Korean BBQ,Happy Lamb Hot Pot"
$ws.Range("B5").Value = "asian_food"
$ws.Range("C5").Value = "additive"

# --- Selection matches target workbook state ---
$ws.Range("A3").Select()
